# Fix the last date entry and append the two missing quarterly dates.
# (commit: "1.fixed some alarm like append()")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last existing row (60) held a wrong date value ("2022-12-19"); correct it.
$ws.Range("A60").Value = "2022-12-31 00:00:00"

# Copy the formatting (text number format / style) of A60 down to the two
# new rows so the new cells share the same style index instead of Excel
# creating a brand-new one.
$ws.Range("A60").Copy()
$ws.Range("A61:A62").PasteSpecial(-4122)

# Append the two new quarterly date strings.
$ws.Range("A61").Value = "2023-03-31 00:00:00"
$ws.Range("A62").Value = "2023-04-30 00:00:00"

# Move the active selection onto the newly appended last cell, matching
# what Excel records after typing/appending the new row.
$ws.Range("A62").Select()

Write-Host "Appended 2023-03-31 and 2023-04-30 rows; fixed 2022-12-31 row."
